$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Data for the new rows: state, election_year, partyname_short, manifesto_url
# (election_manifesto / column D is always TRUE, copied below from an existing
# TRUE cell so it is stored as literal text "TRUE", matching the source data,
# rather than being auto-coerced into an Excel boolean.)
$newRows = @(
    @{ State = "BY"; Year = 2023; Party = "CSU"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/csuregierungsprogramm2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "Grüne"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/gruneregierungsprogramm2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "FW-BY"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fwwahlprogramm2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "AfD"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/afd-bayern.pdf" }
    @{ State = "BY"; Year = 2023; Party = "SPD"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/spdwahlprogramm.pdf" }
    @{ State = "BY"; Year = 2023; Party = "FDP"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fdp-landtagswahlprogramm-2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "Linke"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/linkeoppositionsprogramm2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "BP"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/wahlprogramm-der-bayernpartei-zur-landtagswahl-2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "ÖDP"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/odpwahlprogramm2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "Die Tierschutzpartei"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/wahlflyertierschutzbayern.pdf" }
    @{ State = "BY"; Year = 2023; Party = "Die Humanisten"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/pdhwahlprogramm2023.pdf" }
    @{ State = "BY"; Year = 2023; Party = "dieBasis"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/diebasis-by-wahlprogramm-ltw2023-230105.pdf" }
    @{ State = "HE"; Year = 2023; Party = "CDU"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/cduwahlprogrammltwhessenoffiziell.pdf" }
    @{ State = "HE"; Year = 2023; Party = "Grüne"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/regierungsprogramm-gruene-hessen.pdf" }
    @{ State = "HE"; Year = 2023; Party = "SPD"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/spdhessenwahlprogramm2023v3.pdf" }
    @{ State = "HE"; Year = 2023; Party = "AfD"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/afd-hessen-wahlprogramm-zur-ltw23.pdf" }
    @{ State = "HE"; Year = 2023; Party = "FDP"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fdphe23ltw-programm2023-1.pdf" }
    @{ State = "HE"; Year = 2023; Party = "Linke"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/linkewahlprogramm-ltw2023.pdf" }
    @{ State = "HE"; Year = 2023; Party = "FW-HE"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/fw-wahlprogramm-2023.pdf" }
    @{ State = "HE"; Year = 2023; Party = "DiePartei"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/wpdiepartei2023.pdf" }
    @{ State = "HE"; Year = 2023; Party = "ÖDP"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/oedplandtagswahlhessen2023.pdf" }
    @{ State = "HE"; Year = 2023; Party = "Die Humanisten"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/pdh-hessen-landtagswahl-2023.pdf" }
    @{ State = "HE"; Year = 2023; Party = "dieBasis"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/diebasis-hessen.pdf" }
    @{ State = "HE"; Year = 2023; Party = "DKP"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/dkp-wahlprogramm-2023.pdf" }
    @{ State = "HE"; Year = 2023; Party = "Volt"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/voltlandtagswahlprogrammhessen2023.pdf" }
    @{ State = "HE"; Year = 2023; Party = "WKH"; Url = "https://www.abgeordnetenwatch.de/sites/default/files/election-program-files/klimaliste-hessen-wahlprogramm.pdf" }
)

$startRow = 354
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $rec = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $rec.State
    $ws.Cells.Item($row, 2).Value = $rec.Year
    $ws.Cells.Item($row, 3).Value = $rec.Party
    $ws.Range("D2").Copy($ws.Cells.Item($row, 4))
    $ws.Cells.Item($row, 5).Value = $rec.Url
}

$ws.Application.CutCopyMode = $false

# Leave the view scrolled near the newly added rows, with the same cell
# selected as in the authored workbook (a few rows below the last new row).
$excel.ActiveWindow.ScrollRow = 357
[void]$ws.Range("E383").Select()
